$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the trailing empty paragraph that sits right after
# "Let me know and I'll prepare them for you." (the last paragraph in
# the document body, before sectPr). We index directly rather than via
# .Previous/.Next (unreliable in this host) and rather than assuming a
# fixed paragraph number (robust to any earlier structural differences).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# Give the (currently empty) trailing paragraph a Hindi (India) bidi
# language tag on its paragraph mark - mirrors what Word stamps when a
# user places the caret there and switches on a bidi keyboard layout,
# producing <w:pPr><w:rPr><w:lang w:bidi="hi-IN"/></w:rPr></w:pPr>.
$lastPara.Range.LanguageIDOther = "hi-IN"

# Append a brand-new paragraph after it and fill it in directly via its
# OOXML so the run/break layout matches exactly: a first run with the
# Kubernetes/EKS/GKE sentence, then a second run that starts with a
# manual line break followed by the Terraform sentence - both runs and
# the paragraph mark itself carry the same hi-IN bidi language tag, and
# the paragraph mark also carries the rFonts cs hint.
$lastPara.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $newPara.Range
$insertPoint.Collapse(0)  # wdCollapseEnd

$paraXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:rPr><w:rFonts w:hint="cs"/><w:lang w:bidi="hi-IN"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:bidi="hi-IN"/></w:rPr>' +
  '<w:t>1. Design, deploy, and manage Kubernetes clusters on Amazon EKS and Google Kubernetes Engine (GKE).</w:t>' +
  '</w:r>' +
  '<w:r><w:rPr><w:lang w:bidi="hi-IN"/></w:rPr>' +
  '<w:br/>' +
  '<w:t>2. Implement and maintain infrastructure as code using Terraform.</w:t>' +
  '</w:r>' +
  '</w:p>'

$insertPoint.InsertXML($paraXml)
